# Updated cryptos list on Sun Jul 23 06:35:35 UTC 2023 with GitHub Actions
#
# Refreshes the Coin/Link/Price/Volume(1h) table on the active sheet.
# Column D ("Price") cells are plain text that often *look* numeric
# (e.g. "0.7409", "242.31"); a leading apostrophe forces Excel to keep
# them as text instead of silently coercing them to Double values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.919.49"
$ws.Range("E2").Value = "  +0.07%  "

$ws.Range("D3").Value = "'1.874.87"
$ws.Range("E3").Value = "  -0.97%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'0.7409"
$ws.Range("E5").Value = "  -3.97%  "

$ws.Range("D6").Value = "'242.31"
$ws.Range("E6").Value = "  -0.77%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").Value = "'0.3153"
$ws.Range("E8").Value = "  +0.90%  "

$ws.Range("D9").Value = "'0.07181"
$ws.Range("E9").Value = "  -0.52%  "

$ws.Range("D10").Value = "'24.74"
$ws.Range("E10").Value = "  -3.58%  "

$ws.Range("D11").Value = "'0.08470"
$ws.Range("E11").Value = "  -2.55%  "

$ws.Range("D12").Value = "'0.7533"
$ws.Range("E12").Value = "  -2.22%  "

# Rows 13 & 14 swapped places (WrappedEther <-> Polkadot)
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'5.399"
$ws.Range("E13").Value = "  +0.15%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "'1.871.92"
$ws.Range("E14").Value = "  -6.03%  "

$ws.Range("D15").Value = "'92.61"
$ws.Range("E15").Value = "  -1.70%  "

# Rows 16 & 17 swapped places (Uniswap <-> WrappedBTC)
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "'29.942.13"
$ws.Range("E16").Value = "  -0.73%  "

$ws.Range("B17").Value = "Uniswap"
$ws.Range("C17").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D17").Value = "'6.115"
$ws.Range("E17").Value = "  -1.68%  "

$ws.Range("E18").Value = "  -2.31%  "

$ws.Range("D19").Value = "'243.35"
$ws.Range("E19").Value = "  -0.73%  "

$ws.Range("D20").Value = "'0.000007828"
$ws.Range("E20").Value = "  -0.39%  "

$ws.Range("E21").Value = "  -0.10%  "

$ws.Range("D22").Value = "'2.125.62"
$ws.Range("E22").Value = "  -10.14%  "

$ws.Range("D23").Value = "'7.998"
$ws.Range("E23").Value = "  -2.13%  "

$ws.Range("D24").Value = "'0.9989"
$ws.Range("E24").Value = "  -0.21%  "

$ws.Range("D25").Value = "'0.1557"
$ws.Range("E25").Value = "  -2.30%  "

$ws.Range("D26").Value = "'9.316"
$ws.Range("E26").Value = "  -1.98%  "

$ws.Range("D27").Value = "'165.67"
$ws.Range("E27").Value = "  +2.12%  "

$ws.Range("D28").Value = "'18.63"
$ws.Range("E28").Value = "  -0.87%  "

$ws.Range("D29").Value = "'2.044"
$ws.Range("E29").Value = "  +0.24%  "

$ws.Range("E30").Value = "  +2.79%  "

$ws.Range("D31").Value = "'4.606"
$ws.Range("E31").Value = "  +1.90%  "

$ws.Range("D32").Value = "'1.533"
$ws.Range("E32").Value = "  -0.71%  "

$ws.Range("D33").Value = "'4.283"
$ws.Range("E33").Value = "  +4.10%  "

$ws.Range("D34").Value = "'0.05342"
$ws.Range("E34").Value = "  -2.30%  "

$ws.Range("D35").Value = "'1.245"
$ws.Range("E35").Value = "  -0.21%  "

$ws.Range("D36").Value = "'0.7567"
$ws.Range("E36").Value = "  +0.27%  "

$ws.Range("D37").Value = "'0.9985"
$ws.Range("E37").Value = "  -0.57%  "

$ws.Range("D38").Value = "'2.692"
$ws.Range("E38").Value = "  -0.68%  "

$ws.Range("D39").Value = "'0.01963"
$ws.Range("E39").Value = "  -0.24%  "

$ws.Range("D40").Value = "'2.752"
$ws.Range("E40").Value = "  -1.29%  "

$ws.Range("D41").Value = "'0.4483"
$ws.Range("E41").Value = "  -0.51%  "

$ws.Range("D42").Value = "'1.112.69"
$ws.Range("E42").Value = "  +1.52%  "

$ws.Range("D43").Value = "'6.112"
$ws.Range("E43").Value = "  +1.09%  "

$ws.Range("D44").Value = "'72.47"
$ws.Range("E44").Value = "  -1.90%  "

$ws.Range("D45").Value = "'0.8597"
$ws.Range("E45").Value = "  +0.57%  "

$ws.Range("D46").Value = "'1.001"
$ws.Range("E46").Value = "  +0.11%  "

$ws.Range("D47").Value = "'103.13"
$ws.Range("E47").Value = "  +0.06%  "

$ws.Range("D48").Value = "'7.685"
$ws.Range("E48").Value = "  +0.78%  "

# Rows 49 & 50 swapped places (RenderToken <-> SynthetixNetwork)
$ws.Range("B49").Value = "SynthetixNetwork"
$ws.Range("C49").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D49").Value = "'3.066"
$ws.Range("E49").Value = "  +3.82%  "

$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").Value = "'1.840"
$ws.Range("E50").Value = "  -2.28%  "

$ws.Range("D51").Value = "'2.023.07"
$ws.Range("E51").Value = "  -8.14%  "
